$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Resize the anchored "Frame1" text-box shape (letterhead date/frame box).
#    wp:extent goes 2518410x1798320 EMU -> 2519680x1799590 EMU, i.e.
#    198.3pt x 141.6pt -> 198.4pt x 141.7pt.
# ---------------------------------------------------------------------------
$shape = $d.Shapes.Item(1)
$shape.Width = 198.4
$shape.Height = 141.7

# ---------------------------------------------------------------------------
# 2) Merge the three runs that spelled out the "A-Post-Plus" addressee
#    placeholder into a single run, and swap GRUNDEIGENTUEMER for
#    VERTRETER_NAME_ADDRESS.
# ---------------------------------------------------------------------------
$old = "{{GRUNDEIGENTUEMER or GESUCHSTELLER_NAME_ADDRESS}}(Beilagen gem. Ziff. 2.1 und 2.2) (inkl. Beilagen gem" + [char]0x00E4 + "ss Ziff. (Ziff))"
$new = "{{VERTRETER_NAME_ADDRESS or GESUCHSTELLER_NAME_ADDRESS}}(Beilagen gem. Ziff. 2.1 und 2.2) (inkl. Beilagen gem" + [char]0x00E4 + "ss Ziff. (Ziff))"
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Register the additional list-label character styles (ListLabel75 ..
#    ListLabel92) that back the numbering fonts used further down the
#    numbering definitions.
# ---------------------------------------------------------------------------
$labelFonts = @(
  "Arial",       # 75
  "Courier New", # 76
  "Wingdings",   # 77
  "Symbol",      # 78
  "Courier New", # 79
  "Wingdings",   # 80
  "Symbol",      # 81
  "Courier New", # 82
  "Wingdings",   # 83
  "Arial",       # 84
  "Courier New", # 85
  "Wingdings",   # 86
  "Symbol",      # 87
  "Courier New", # 88
  "Wingdings",   # 89
  "Symbol",      # 90
  "Courier New", # 91
  "Wingdings"    # 92
)

for ($i = 0; $i -lt $labelFonts.Count; $i++) {
  $num = 75 + $i
  $style = $d.Styles.Add("ListLabel " + $num, 2)
  $style.QuickStyle = $true
  $style.Font.NameBi = $labelFonts[$i]
}
